$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 217.75
$ws.Range("I8").Value = 217.75
$ws.Range("K8").Value = 653.25
$ws.Range("M8").Value = -514.25
$ws.Range("H18").Value = 1000
$ws.Range("I18").Value = 1000
$ws.Range("K18").Value = 1000
$ws.Range("M18").Value = -716
$ws.Range("H20").Value = 7299.3335
$ws.Range("I20").Value = 7299.3335
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 7299.3335
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -7069.3335
$ws.Range("N20").ClearContents()
$ws.Range("H35").Value = 7299.3335
$ws.Range("I35").Value = 7299.3335
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 7299.3335
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -6920.3335
$ws.Range("N35").ClearContents()
$ws.Range("H51").Value = 14999.25
$ws.Range("J51").Value = 14999.25
$ws.Range("L51").Value = 14999.25
$ws.Range("N51").Value = -15967.25
$ws.Range("H74").Value = 6141.0557
$ws.Range("I74").Value = 6141.0557
$ws.Range("K74").Value = 6141.0557
$ws.Range("M74").Value = -5205.0557
$ws.Range("H77").Value = 6141.0557
$ws.Range("I77").Value = 6141.0557
$ws.Range("K77").Value = 30705.2785
$ws.Range("M77").Value = -26025.2785
$ws.Range("H86").Value = 4673.75
$ws.Range("I86").Value = 4484.2856
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 4484.2856
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -3361.2856
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 4673.75
$ws.Range("I89").Value = 4484.2856
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 22421.428
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -16805.428
$ws.Range("N89").Value = -41232
$ws.Range("H98").Value = 967.8461
$ws.Range("I98").Value = 967.8461
$ws.Range("K98").Value = 967.8461
$ws.Range("M98").Value = 530.1539
$ws.Range("H115").Value = 1955.3334
$ws.Range("J115").Value = 10000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -33134
$ws.Range("H122").Value = 967.8461
$ws.Range("I122").Value = 967.8461
$ws.Range("K122").Value = 2903.5383
$ws.Range("M122").Value = -453.5383000000002
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H138").Value = 3569.375
$ws.Range("I138").Value = 3823.7856
$ws.Range("J138").Value = 3371.5
$ws.Range("K138").Value = 11471.3568
$ws.Range("L138").Value = 10114.5
$ws.Range("M138").Value = -6331.356800000001
$ws.Range("N138").Value = -20394.5
$ws.Range("H140").Value = 74000
$ws.Range("J140").Value = 74000
$ws.Range("L140").Value = 74000
$ws.Range("N140").Value = -84360
# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1243.9048
$ws.Range("I2").Value = 810.2143
$ws.Range("K2").Value = 810.2143
$ws.Range("M2").Value = -697.2143
$ws.Range("H32").Value = 9491.842000000001
$ws.Range("I32").Value = 6023.2666
$ws.Range("K32").Value = 6023.2666
$ws.Range("M32").Value = -5736.2666
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 0
$ws.Range("K33").Value = 0
$ws.Range("M33").ClearContents()
$ws.Range("H61").Value = 2192.5
$ws.Range("I61").Value = 2090.1667
$ws.Range("K61").Value = 2090.1667
$ws.Range("M61").Value = -1878.1667
$ws.Range("H102").Value = 2767.5293
$ws.Range("I102").Value = 2503.5
$ws.Range("K102").Value = 2503.5
$ws.Range("M102").Value = -881.5
$ws.Range("H116").Value = 1243.9048
$ws.Range("I116").Value = 810.2143
$ws.Range("K116").Value = 810.2143
$ws.Range("M116").Value = 1483.7857
$ws.Range("H122").Value = 2672.48
$ws.Range("I122").Value = 2810.3044
$ws.Range("K122").Value = 8430.913199999999
$ws.Range("M122").Value = -5980.913199999999
$ws.Range("H136").Value = 2192.5
$ws.Range("I136").Value = 2090.1667
$ws.Range("K136").Value = 6270.500100000001
$ws.Range("M136").Value = -3720.500100000001
# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1243.9048
$ws.Range("I3").Value = 810.2143
$ws.Range("K3").Value = 810.2143
$ws.Range("M3").Value = -696.2143
$ws.Range("H94").Value = 1021.8461
$ws.Range("I94").Value = 1112.6666
$ws.Range("J94").Value = 944
$ws.Range("K94").Value = 1112.6666
$ws.Range("L94").Value = 944
$ws.Range("M94").Value = -661.6666
$ws.Range("N94").Value = -1846
# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H86").Value = 13413.929
$ws.Range("I86").Value = 11445.444
$ws.Range("J86").Value = 16957.2
$ws.Range("K86").Value = 11445.444
$ws.Range("L86").Value = 16957.2
$ws.Range("M86").Value = -10322.444
$ws.Range("N86").Value = -19203.2
$ws.Range("H89").Value = 13413.929
$ws.Range("I89").Value = 11445.444
$ws.Range("J89").Value = 16957.2
$ws.Range("K89").Value = 57227.22
$ws.Range("L89").Value = 84786
$ws.Range("M89").Value = -51611.22
$ws.Range("N89").Value = -96018
# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H24").Value = 1851.8
$ws.Range("I24").Value = 1440
$ws.Range("J24").Value = 3499
$ws.Range("K24").Value = 4320
$ws.Range("L24").Value = 10497
$ws.Range("M24").Value = -4090
$ws.Range("N24").Value = -10957
$ws.Range("H80").Value = 4993.5
$ws.Range("I80").Value = 4989
$ws.Range("J80").Value = 4995
$ws.Range("K80").Value = 14967
$ws.Range("L80").Value = 14985
$ws.Range("M80").Value = -14031
$ws.Range("N80").Value = -16857
$ws.Range("H83").Value = 4993.5
$ws.Range("I83").Value = 4989
$ws.Range("J83").Value = 4995
$ws.Range("K83").Value = 44901
$ws.Range("L83").Value = 44955
$ws.Range("M83").Value = -40221
$ws.Range("N83").Value = -54315
$ws.Range("H107").Value = 450.18518
$ws.Range("J107").Value = 451.6875
$ws.Range("L107").Value = 1355.0625
$ws.Range("N107").Value = -5195.0625
$ws.Range("H131").Value = 2609.6365
$ws.Range("J131").Value = 4423.5
$ws.Range("L131").Value = 13270.5
$ws.Range("N131").Value = -23350.5
# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 9899.5
$ws.Range("J24").Value = 14900
$ws.Range("L24").Value = 14900
$ws.Range("N24").Value = -15246
$ws.Range("H36").Value = 2750
$ws.Range("I36").Value = 1500
$ws.Range("K36").Value = 1500
$ws.Range("M36").Value = -1015
$ws.Range("H48").Value = 5000
$ws.Range("I48").Value = 5000
$ws.Range("K48").Value = 5000
$ws.Range("M48").Value = -4515
$ws.Range("H70").Value = 2465
$ws.Range("I70").Value = 2465
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 2465
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -2195
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 2465
$ws.Range("I73").Value = 2465
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 2465
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -1529
$ws.Range("N73").ClearContents()
$ws.Range("H92").Value = 12523.667
$ws.Range("J92").Value = 10260.5
$ws.Range("L92").Value = 10260.5
$ws.Range("N92").Value = -14004.5
$ws.Range("H102").Value = 1648.6666
$ws.Range("I102").Value = 1648.6666
$ws.Range("K102").Value = 1648.6666
$ws.Range("M102").Value = -26.66660000000002
$ws.Range("H132").Value = 14496504
$ws.Range("I132").Value = 3911
$ws.Range("K132").Value = 11733
$ws.Range("M132").Value = -9203
$ws.Range("H136").Value = 65996.14
$ws.Range("J136").Value = 65996.14
$ws.Range("L136").Value = 197988.42
$ws.Range("N136").Value = -203088.42
# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 83335810
$ws.Range("J22").Value = 125003110
$ws.Range("L22").Value = 125003110
$ws.Range("N22").Value = -125003700
$ws.Range("H27").Value = 83335810
$ws.Range("J27").Value = 125003110
$ws.Range("L27").Value = 125003110
$ws.Range("N27").Value = -125003324
$ws.Range("H46").Value = 6914.4546
$ws.Range("I46").Value = 6222.7144
$ws.Range("K46").Value = 6222.7144
$ws.Range("M46").Value = -6034.7144
$ws.Range("H132").Value = 3010.889
$ws.Range("J132").Value = 2626.5334
$ws.Range("L132").Value = 7879.600199999999
$ws.Range("N132").Value = -12939.6002
# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 40999.8
$ws.Range("J46").Value = 40999.8
$ws.Range("L46").Value = 40999.8
$ws.Range("N46").Value = -41461.8
$ws.Range("H122").Value = 2399.25
$ws.Range("I122").Value = 2353
$ws.Range("K122").Value = 7059
$ws.Range("M122").Value = -4609
$ws.Range("H134").Value = 40999.8
$ws.Range("J134").Value = 40999.8
$ws.Range("L134").Value = 122999.4
$ws.Range("N134").Value = -128069.4
$ws.Range("H136").Value = 7765.647
$ws.Range("I136").Value = 8192
$ws.Range("J136").Value = 944
$ws.Range("K136").Value = 24576
$ws.Range("L136").Value = 2832
$ws.Range("M136").Value = -22026
$ws.Range("N136").Value = -7932
